# Brief 5 Requirements (1)
# Mark the Week 5 "Required" checklist tasks (rows 35-38, 41, 47, 48) as
# completed. These cells are linked to Form-Control checkboxes
# (Check Box 46, 47, 48, 49, 52, 56, 57); setting the linked cell to TRUE
# both "checks" the box and drives the dependent IF()/COUNTIFS() formulas
# (E35:E38/E41/E47/E48 -> "Done", N7/N8/N10/N11 summary stats) to
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$doneCells = @("I35", "I36", "I37", "I38", "I41", "I47", "I48")
foreach ($cellRef in $doneCells) {
    $ws.Range($cellRef).Value = $true
}
